$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = -0.000000000008779167981424272
$ws.Range("B2").Value = -0.00000000001609528544845711
$ws.Range("C2").Value = -0.00000000009657171269114809
$ws.Range("D2").Value = 0.00000000000000000002984846159787282
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.0000000000000000000006217560184328447
$ws.Range("G2").Value = 0.0000000000000000000001727411654803868
$ws.Range("H2").Value = 0.00000000000000000000000000004692058802350669
$ws.Range("A3").Value = -0.0000000000000001478772354528204
$ws.Range("B3").Value = -0.0000000000000004318785482454018
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 0.0000000003380523074071533
$ws.Range("E3").Value = -0.00000000009657171269114488
$ws.Range("F3").Value = -0.000000000019314236368038
$ws.Range("G3").Value = -0.000000000009656633869885605
$ws.Range("H3").Value = 0.0000000000000000001799092274089556
$ws.Range("A4").Value = -0.000000000000001154681825411029
$ws.Range("B4").Value = -1
$ws.Range("C4").Value = -0.0000000000000001357950826736785
$ws.Range("D4").Value = -0.00000000001931434253822938
$ws.Range("E4").Value = -0.00000000001609528544882413
$ws.Range("F4").Value = -0.00000000000000005875310241946868
$ws.Range("G4").Value = -0.000000000006438105270984324
$ws.Range("H4").Value = -0.00000000000000000000002061768773616638
$ws.Range("A5").Value = -0.00000000000000006995381887074578
$ws.Range("B5").Value = -0.000000000006437929010255724
$ws.Range("C5").Value = -0.00000000000965763618108535
$ws.Range("D5").Value = 0.000000000000000000003108536132856199
$ws.Range("E5").Value = -0.000000000000000000001208949519111344
$ws.Range("F5").Value = 0.000000000000000008541667303886157
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.000000000004598652985294717
$ws.Range("A6").Value = -1
$ws.Range("B6").Value = 0.0000000000000008405394991419988
$ws.Range("C6").Value = -0.0000000000000002092220421464854
$ws.Range("D6").Value = -0.000000000009657171269114646
$ws.Range("E6").Value = -0.000000000008779167983118333
$ws.Range("F6").Value = -0.000000000006438084802858598
$ws.Range("G6").Value = -0.000000000000000002547768570861929
$ws.Range("H6").Value = -0.00000000000000000000001612184258303316
$ws.Range("A7").Value = -0.000000000006438114179409815
$ws.Range("B7").Value = 0.000000000000000000000000003090924084001512
$ws.Range("C7").Value = -0.00000000001931461135188511
$ws.Range("D7").Value = 0.000000000000000000006402104597874661
$ws.Range("E7").Value = -0.000000000000000000002543480137108714
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.0000000000000001942890383709254
$ws.Range("H7").Value = 0.000000000003714296641956536
$ws.Range("A8").Value = -0.000000000009657171269114728
$ws.Range("B8").Value = -0.0000000000193143425382294
$ws.Range("C8").Value = 0.00000000002897151380734415
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.00000000000000000000240217692243337
$ws.Range("F8").Value = 0.00000000000000000000005246550266921314
$ws.Range("G8").Value = 0.00000000000000000000001780216013899641
$ws.Range("H8").Value = 0.000000000003115111467656421
$ws.Range("A9").Value = -0.00000000000000000000003791118289460318
$ws.Range("B9").Value = -0.00000000000000000000006908155685931354
$ws.Range("C9").Value = -0.0000000000000000008419587640931153
$ws.Range("D9").Value = 0.000000000003115352420863995
$ws.Range("E9").Value = -0.000000000000000000000000000000004346547728067738
$ws.Range("F9").Value = 0.000000000003714296641936368
$ws.Range("G9").Value = 0.000000000004598652985297356
$ws.Range("H9").Value = -1
